$d = $word.ActiveDocument

# Ordered list of (old, new) text replacements. The order is chosen so that
# at the moment each replacement runs, the "old" string occurs exactly once in
# the document (this matters because "73÷6=" is both a target value for one cell
# and the newly-introduced value for another cell).
$replacements = @(
    @("87÷5=", "33÷6="),
    @("62÷7=", "23÷3="),
    @("73÷6=", "60÷6="),
    @("74÷5=", "73÷6="),
    @("14÷5=", "43÷2="),
    @("35÷7=", "80÷3="),
    @("83÷3=", "77÷6="),
    @("62÷2=", "67÷6="),
    @("37÷8=", "92÷8="),
    @("49÷7=", "13÷2="),
    @("86÷6=", "47÷3="),
    @("16÷4=", "81÷3="),
    @("96÷9=", "55÷7="),
    @("98÷7=", "55÷4="),
    @("33÷7=", "41÷9="),
    @("59÷6=", "13÷6="),
    @("10÷8=", "82÷6="),
    @("86÷4=", "65÷5="),
    @("94÷2=", "49÷6="),
    @("26÷7=", "29÷3="),
    @("54÷8=", "32÷2="),
    @("30÷5=", "83÷9="),
    @("50÷4=", "68÷3="),
    @("67÷3=", "35÷2="),
    @("57÷8=", "62÷8="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $found) {
        Write-Host "WARNING: could not find text to replace: $old"
    }
}

$d.Save()
